$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date formatting from the last existing date cell (A7) onto the
# two new date cells so they pick up the same style index (short date fmt)
# rather than Excel inventing a brand new number format.
$ws.Range("A7").Copy()
$ws.Range("A8:A9").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A8").Value = 45969
$ws.Range("B8").Value = 2

$ws.Range("A9").Value = 45968
$ws.Range("B9").Value = 3

$ws.Range("A8:B9").Select()
